# Updated cryptos list with refreshed Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.129.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "'3.138.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'535.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "

$ws.Range("D6").Value = "'139.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.505"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.17%  "

$ws.Range("D9").Value = "'7.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("E10").Value = "  +2.33%  "

$ws.Range("D11").Value = "'0.418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.28%  "

$ws.Range("E12").Value = "  +3.19%  "

$ws.Range("D13").Value = "'3.678.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "'25.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  +5.84%  "

$ws.Range("D16").Value = "'58.212.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D18").Value = "'3.137.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").Value = "'12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.57%  "

$ws.Range("E20").Value = "  +3.99%  "

$ws.Range("D21").Value = "'375.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.29%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "'5.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "

$ws.Range("D24").Value = "'70.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").Value = "'0.0₃0887"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").Value = "'7.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.82%  "

$ws.Range("D30").Value = "'6.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.19%  "

$ws.Range("D31").Value = "'1.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("E32").Value = "  +4.04%  "

$ws.Range("E33").Value = "  +5.99%  "

$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").Value = "'161.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("E36").Value = "  +3.54%  "

$ws.Range("D37").Value = "'1.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.03%  "

$ws.Range("D38").Value = "'25.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  +6.03%  "

$ws.Range("D40").Value = "'2.620.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.46%  "

$ws.Range("D41").Value = "'0.0674"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("E42").Value = "  +4.05%  "

$ws.Range("D43").Value = "'39.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.28%  "

$ws.Range("D44").Value = "'0.701"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "

$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("E47").Value = "  +4.38%  "

$ws.Range("E48").Value = "  +2.83%  "

$ws.Range("D49").Value = "'0.0994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.47%  "

$ws.Range("E50").Value = "  +3.03%  "

$ws.Range("D51").Value = "'0.752"
$ws.Range("D51").Style = "Normal"
